$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.002.12"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "2.295.01"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'300.00"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").Value = "'97.69"
$ws.Range("E6").Value = "  -2.87%  "

$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = "  +2.43%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'36.05"
$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "'17.63"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  -1.27%  "

$ws.Range("D15").Value = "2.652.49"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").Value = "2.294.82"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("E17").Value = "  -1.52%  "

$ws.Range("D18").Value = "42.926.72"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").Value = "'12.75"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +0.97%  "

$ws.Range("D21").Value = "'6.11"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "'68.88"
$ws.Range("E22").Value = "  +1.35%  "

$ws.Range("D23").Value = "'236.92"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  -1.05%  "

$ws.Range("D27").Value = "'24.91"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "'165.14"
$ws.Range("E28").Value = "  -1.57%  "

$ws.Range("D29").Value = "'2.04"
$ws.Range("E29").Value = "  -0.56%  "

$ws.Range("D30").Value = "'9.07"
$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("D31").Value = "'33.04"
$ws.Range("E31").Value = "  -4.46%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  +0.35%  "

$ws.Range("D34").Value = "'4.73"
$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("D35").Value = "'17.80"
$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("E41").Value = "  -1.85%  "

$ws.Range("D42").Value = "2.004.82"
$ws.Range("E42").Value = "  +1.13%  "

$ws.Range("E43").Value = "  -2.52%  "

$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("D45").Value = "'10.23"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").Value = "'17.41"
$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("E47").Value = "  -2.95%  "

$ws.Range("D48").Value = "'54.06"
$ws.Range("E48").Value = "  -2.36%  "

$ws.Range("D49").Value = "2.519.61"
$ws.Range("E49").Value = "  -0.46%  "

$ws.Range("D50").Value = "'73.22"
$ws.Range("E50").Value = "  +3.34%  "

$ws.Range("E51").Value = "  -1.86%  "
